$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 436, shifting existing rows 436:536 down to 437:537.
$ws.Rows.Item(436).Insert()

# Populate the newly inserted row with the new data record.
$ws.Cells.Item(436, 1).Value = 3
$ws.Cells.Item(436, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(436, 3).Value = "Coquimbo"
$ws.Cells.Item(436, 4).Value = 44785
$ws.Cells.Item(436, 5).Value = 5
$ws.Cells.Item(436, 6).Value = 100112037
$ws.Cells.Item(436, 7).Value = "Cebollín"
$ws.Cells.Item(436, 8).Value = "Sin especificar"
$ws.Cells.Item(436, 9).Value = "Primera"
$ws.Cells.Item(436, 10).Value = 290
$ws.Cells.Item(436, 11).Value = 6000
$ws.Cells.Item(436, 12).Value = 7000
$ws.Cells.Item(436, 13).Value = 6448
$ws.Cells.Item(436, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(436, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(436, 16).Value = 179
$ws.Cells.Item(436, 17).Value = 36
$ws.Cells.Item(436, 18).Value = "Hortaliza"
